$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column mapping: A=1, B=2, C=3, D=4, E=5

$ws.Cells.Item(3, 4).Value = -7.811999999999999
$ws.Cells.Item(4, 4).Value = -8.045

$ws.Cells.Item(6, 5).Value = 12.718

$ws.Cells.Item(7, 4).Value = -8.102
$ws.Cells.Item(7, 5).Value = 12.895

$ws.Cells.Item(8, 4).Value = -8.010999999999999
$ws.Cells.Item(8, 5).Value = 13.03

$ws.Cells.Item(11, 1).Value = -21.556

$ws.Cells.Item(12, 1).Value = -21.952
$ws.Cells.Item(12, 4).Value = -8.029999999999999

$ws.Cells.Item(14, 4).Value = -8.224000000000002

$ws.Cells.Item(15, 1).Value = -21.221

$ws.Cells.Item(19, 5).Value = 12.169

$ws.Cells.Item(21, 5).Value = 13.136

$ws.Cells.Item(22, 4).Value = -7.812

$ws.Cells.Item(24, 5).Value = 12.623

$ws.Cells.Item(25, 5).Value = 12.169
